$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) New "consignes_TP6" entry: insert a fresh row at row 6 (pushes the
#    "architecture" block, and everything below, down by one row).
$ws.Rows(6).Insert()
$ws.Range("A6").Value = "consignes"
$ws.Range("B6").Value = "consignes_TP6"
$ws.Range("C6").Value = "Consignes"

# 2) New "fichiers" domain rows (TP06 subject), inserted right before the
#    trailing "fin" marker row (which is now row 70 after the insert above).
$ws.Rows(70).Insert()
$ws.Rows(70).Insert()
$ws.Rows(70).Insert()
$ws.Rows(70).Insert()
$ws.Rows(70).Insert()

$ws.Range("A70").Value = "fichiers"
$ws.Range("B70").Value = "FIC-007"
$ws.Range("C70").Value = 'Classement de la "Saintélyon"'

$ws.Range("A71").Value = "fichiers"
$ws.Range("B71").Value = "FIC-009"
$ws.Range("C71").Value = "Analyse d'un dipôle électrique"

$ws.Range("A72").Value = "fichiers"
$ws.Range("B72").Value = "FIC-010"
$ws.Range("C72").Value = "Traitement de données physiologiques"

$ws.Range("A73").Value = "fichiers"
$ws.Range("B73").Value = "FIC-011"
$ws.Range("C73").Value = "Lecture d'un texte"

$ws.Range("A74").Value = "fichiers"
$ws.Range("B74").Value = "FIC-012"
$ws.Range("C74").Value = "Résultats de l'Embrunman"

# 3) Leave the selection where the author left it at commit time.
$ws.Range("C74").Select()
